# Auto-generated from the OOXML diff: updates currentAveragePrice /
# LevePrice / LeveProfit derived columns (H,I,J,K,L,M,N) for the affected
# rows across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9000531
$ws.Range("J17").Value = 9000531
$ws.Range("L17").Value = 27001593
$ws.Range("N17").Value = -27001929
$ws.Range("H53").Value = 434.6154
$ws.Range("J53").Value = 532.4
$ws.Range("L53").Value = 532.4
$ws.Range("N53").Value = -1806.4
$ws.Range("H75").Value = 47216.6
$ws.Range("J75").Value = 47216.6
$ws.Range("L75").Value = 47216.6
$ws.Range("N75").Value = -49088.6
$ws.Range("H76").Value = 7498.6665
$ws.Range("I76").Value = 7122.5
$ws.Range("K76").Value = 7122.5
$ws.Range("M76").Value = -6807.5
$ws.Range("H78").Value = 47216.6
$ws.Range("J78").Value = 47216.6
$ws.Range("L78").Value = 141649.8
$ws.Range("N78").Value = -151009.8
$ws.Range("H79").Value = 7498.6665
$ws.Range("I79").Value = 7122.5
$ws.Range("K79").Value = 7122.5
$ws.Range("M79").Value = -6030.5
$ws.Range("H112").Value = 1859.2858
$ws.Range("I112").Value = 992.5
$ws.Range("J112").Value = 2647.2727
$ws.Range("K112").Value = 2977.5
$ws.Range("L112").Value = 7941.8181
$ws.Range("M112").Value = -1869.5
$ws.Range("N112").Value = -10157.8181
$ws.Range("H138").Value = 9843.421
$ws.Range("I138").Value = 3099.4
$ws.Range("J138").Value = 12252
$ws.Range("K138").Value = 9298.200000000001
$ws.Range("L138").Value = 36756
$ws.Range("M138").Value = -4158.200000000001
$ws.Range("N138").Value = -47036

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2287.889
$ws.Range("I45").Value = 1620.5
$ws.Range("J45").Value = 2821.8
$ws.Range("K45").Value = 1620.5
$ws.Range("L45").Value = 2821.8
$ws.Range("M45").Value = -1243.5
$ws.Range("N45").Value = -3575.8
$ws.Range("H61").Value = 15629057
$ws.Range("I61").Value = 15629057
$ws.Range("K61").Value = 15629057
$ws.Range("M61").Value = -15628845
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 21740960
$ws.Range("I74").Value = 30304426
$ws.Range("J74").Value = 2931.923
$ws.Range("K74").Value = 30304426
$ws.Range("L74").Value = 2931.923
$ws.Range("M74").Value = -30303552
$ws.Range("N74").Value = -4679.923
$ws.Range("H76").Value = 48749.5
$ws.Range("J76").Value = 48749.5
$ws.Range("L76").Value = 48749.5
$ws.Range("N76").Value = -49425.5
$ws.Range("H77").Value = 21740960
$ws.Range("I77").Value = 30304426
$ws.Range("J77").Value = 2931.923
$ws.Range("K77").Value = 151522130
$ws.Range("L77").Value = 14659.615
$ws.Range("M77").Value = -151517762
$ws.Range("N77").Value = -23395.615
$ws.Range("H79").Value = 48749.5
$ws.Range("J79").Value = 48749.5
$ws.Range("L79").Value = 48749.5
$ws.Range("N79").Value = -51089.5
$ws.Range("H97").Value = 1422.88
$ws.Range("I97").Value = 1162.3636
$ws.Range("J97").Value = 3333.3333
$ws.Range("K97").Value = 1162.3636
$ws.Range("L97").Value = 3333.3333
$ws.Range("M97").Value = -666.3635999999999
$ws.Range("N97").Value = -4325.3333
$ws.Range("H122").Value = 2499.6553
$ws.Range("I122").Value = 1688.1177
$ws.Range("K122").Value = 5064.3531
$ws.Range("M122").Value = -2614.3531
$ws.Range("H132").Value = 90946890
$ws.Range("I132").Value = 20454.166
$ws.Range("K132").Value = 61362.49800000001
$ws.Range("M132").Value = -58832.49800000001
$ws.Range("H136").Value = 15629057
$ws.Range("I136").Value = 15629057
$ws.Range("K136").Value = 46887171
$ws.Range("M136").Value = -46884621
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 76143.14
$ws.Range("J76").Value = 76143.14
$ws.Range("L76").Value = 76143.14
$ws.Range("N76").Value = -76773.14
$ws.Range("H79").Value = 76143.14
$ws.Range("J79").Value = 76143.14
$ws.Range("L79").Value = 76143.14
$ws.Range("N79").Value = -78327.14
$ws.Range("H86").Value = 13383.363
$ws.Range("I86").Value = 14152.375
$ws.Range("K86").Value = 14152.375
$ws.Range("M86").Value = -13029.375
$ws.Range("H89").Value = 13383.363
$ws.Range("I89").Value = 14152.375
$ws.Range("K89").Value = 70761.875
$ws.Range("M89").Value = -65145.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3985.6667
$ws.Range("J31").Value = 5645.857
$ws.Range("L31").Value = 5645.857
$ws.Range("N31").Value = -6235.857
$ws.Range("H34").Value = 3985.6667
$ws.Range("J34").Value = 5645.857
$ws.Range("L34").Value = 5645.857
$ws.Range("N34").Value = -6049.857
$ws.Range("H92").Value = 64999.75
$ws.Range("J92").Value = 64999.75
$ws.Range("L92").Value = 64999.75
$ws.Range("N92").Value = -69991.75
$ws.Range("H132").Value = 2683.6155
$ws.Range("I132").Value = 2365.5833
$ws.Range("K132").Value = 7096.749899999999
$ws.Range("M132").Value = -4566.749899999999
$ws.Range("H133").Value = 54125.2
$ws.Range("J133").Value = 60156.5
$ws.Range("L133").Value = 60156.5
$ws.Range("N133").Value = -65216.5
$ws.Range("H134").Value = 2469.0715
$ws.Range("J134").Value = 5799
$ws.Range("L134").Value = 17397
$ws.Range("N134").Value = -22467

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 558.1
$ws.Range("J98").Value = 618.6
$ws.Range("L98").Value = 1855.8
$ws.Range("N98").Value = -4851.8
$ws.Range("H109").Value = 4648.3
$ws.Range("I109").Value = 4831.4443
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 14494.3329
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = -13454.3329
$ws.Range("N109").Value = -11080
$ws.Range("H137").Value = 3033
$ws.Range("J137").Value = 3033
$ws.Range("L137").Value = 9099
$ws.Range("N137").Value = -19299
$ws.Range("H140").Value = 3163.4167
$ws.Range("I140").Value = 3163.4167
$ws.Range("K140").Value = 9490.250100000001
$ws.Range("M140").Value = -4310.250100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 563.6177
$ws.Range("I97").Value = 570.74194
$ws.Range("J97").Value = 490
$ws.Range("K97").Value = 570.74194
$ws.Range("L97").Value = 490
$ws.Range("M97").Value = -74.74194
$ws.Range("N97").Value = -1482

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 12502.5
$ws.Range("I18").Value = 12502.5
$ws.Range("K18").Value = 12502.5
$ws.Range("M18").Value = -12330.5
$ws.Range("H82").Value = 2023.6923
$ws.Range("I82").Value = 1849.8182
$ws.Range("J82").Value = 2980
$ws.Range("K82").Value = 1849.8182
$ws.Range("L82").Value = 2980
$ws.Range("M82").Value = -1488.8182
$ws.Range("N82").Value = -3702
$ws.Range("H85").Value = 2023.6923
$ws.Range("I85").Value = 1849.8182
$ws.Range("J85").Value = 2980
$ws.Range("K85").Value = 1849.8182
$ws.Range("L85").Value = 2980
$ws.Range("M85").Value = -601.8181999999999
$ws.Range("N85").Value = -5476
$ws.Range("H136").Value = 1430486.4
$ws.Range("I136").Value = 1668167.5
$ws.Range("K136").Value = 5004502.5
$ws.Range("M136").Value = -5001952.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 21624.8
$ws.Range("J74").Value = 21624.8
$ws.Range("L74").Value = 21624.8
$ws.Range("N74").Value = -23496.8
$ws.Range("H77").Value = 21624.8
$ws.Range("J77").Value = 21624.8
$ws.Range("L77").Value = 64874.39999999999
$ws.Range("N77").Value = -74234.39999999999
$ws.Range("H111").Value = 54850
$ws.Range("J111").Value = 54850
$ws.Range("L111").Value = 54850
$ws.Range("N111").Value = -63030
$ws.Range("H132").Value = 3412.2666
$ws.Range("I132").Value = 2986.2727
$ws.Range("J132").Value = 4583.75
$ws.Range("K132").Value = 8958.8181
$ws.Range("L132").Value = 13751.25
$ws.Range("M132").Value = -6428.8181
$ws.Range("N132").Value = -18811.25
$ws.Range("H140").Value = 97143
$ws.Range("J140").Value = 97143
$ws.Range("L140").Value = 97143
$ws.Range("N140").Value = -107503
$ws.Range("H141").Value = 97470
$ws.Range("J141").Value = 97470
$ws.Range("L141").Value = 97470

